$wb = $excel.ActiveWorkbook

# Mapping of worksheet name -> { cell -> new value }
# Derived from the target OOXML diff: MP time-limit change (most "MP solve time" / solve-time
# cells) plus a fix to the fixed-recourse objective data ("objective" column C on Sheet1 and
# one corrected B-column value that also appears on sheet "6").

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("C2").Value = 1252.369127951
$ws.Range("C3").Value = 1404.700269338
$ws.Range("C4").Value = 1196.390677708
$ws.Range("C5").Value = 804.675731879
$ws.Range("C6").Value = 1382.103584759
$ws.Range("B7").Value = -850.8159243656047
$ws.Range("C7").Value = 1373.464032942
$ws.Range("C8").Value = 1192.674797786
$ws.Range("C9").Value = 1132.823321609
$ws.Range("C10").Value = 1495.255522237
$ws.Range("C11").Value = 658.308469182

$ws = $wb.Worksheets.Item("1")
$ws.Range("D2").Value = 12.316470153888794
$ws.Range("C3").Value = 0.09976696032066357
$ws.Range("D3").Value = 1065.0611238394972

$ws = $wb.Worksheets.Item("2")
$ws.Range("D2").Value = 12.578151509451171
$ws.Range("C3").Value = 0.09301862524417995
$ws.Range("D3").Value = 1220.9666475295562

$ws = $wb.Worksheets.Item("3")
$ws.Range("D2").Value = 11.435332970393311
$ws.Range("D3").Value = 989.0963921435509

$ws = $wb.Worksheets.Item("4")
$ws.Range("D2").Value = 11.01385122432251
$ws.Range("D3").Value = 612.0706227661673

$ws = $wb.Worksheets.Item("5")
$ws.Range("D2").Value = 11.44146026940747
$ws.Range("C3").Value = 0.08568619061343723
$ws.Range("D3").Value = 1201.2744398542748

$ws = $wb.Worksheets.Item("6")
$ws.Range("D2").Value = 11.792407551202027
$ws.Range("B3").Value = -850.8159243656047
$ws.Range("C3").Value = 0.09935042989076803
$ws.Range("D3").Value = 1187.7891846597388

$ws = $wb.Worksheets.Item("7")
$ws.Range("D2").Value = 10.902888867528686
$ws.Range("D3").Value = 1009.2711717171753

$ws = $wb.Worksheets.Item("8")
$ws.Range("D2").Value = 11.087962900236938
$ws.Range("D3").Value = 943.7851584157027

$ws = $wb.Worksheets.Item("9")
$ws.Range("D2").Value = 11.842163594135377
$ws.Range("C3").Value = 0.0948664709994042
$ws.Range("D3").Value = 1313.1303340175432

$ws = $wb.Worksheets.Item("10")
$ws.Range("D2").Value = 11.803123490339722
$ws.Range("D3").Value = 474.58255641225463
